$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shp = $s.Shapes.Item(2)
$shp.Table.ApplyStyle("{66C03912-E092-4E07-A942-F72B50BC6314}")
